$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.684.07'
$ws.Range('D3').Value = '1.800.49'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.557'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.34%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.83'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0683'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0935'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '2.059.63'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.71%  '
$ws.Range('D14').Value = '1.798.19'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.644'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').Value = '34.665.86'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '256.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('D20').Value = '0.0₃0785'
$ws.Range('E20').Value = '  +3.50%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.48'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('E24').Value = '  -4.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  -3.03%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.82'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0519'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.44%  '
$ws.Range('D35').Value = '1.459.34'
$ws.Range('E35').Value = '  -4.60%  '
$ws.Range('E36').Value = '  -1.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.638'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '83.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.906'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0507'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.94'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('D46').Value = '1.958.79'
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('E47').Value = '  -3.33%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.98'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '101.40'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '49.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.94%  '
